# This workbook holds quarterly survey values ("Trimestre" / "Valor") in three
# stacked 21-row blocks (Brasil: rows 2-22, Nordeste: rows 23-43, Sergipe: rows
# 44-64). The upload replaces the data with one quarter advanced: every row's
# Trimestre (C) and Valor (D) become what used to be the NEXT row's values,
# and the final row of each block picks up a brand-new quarter that wasn't in
# the sheet before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each block's first data row, last data row, and the brand-new trailing
# quarter's date/value (one quarter after the old last date in that block).
$blocks = @(
    @{ start = 2;  end = 22; newDate = "01/07/2024"; newValue = 93.63718985731164 },
    @{ start = 23; end = 43; newDate = "01/07/2024"; newValue = 91.32792385113521 },
    @{ start = 44; end = 64; newDate = "01/07/2024"; newValue = 91.61462979482606 }
)

foreach ($block in $blocks) {
    $start = $block.start
    $end = $block.end

    # Capture this block's original Trimestre/Valor columns before overwriting
    # any of them (row r reads row r+1's ORIGINAL values, so we must snapshot
    # first rather than mutate in place while reading ahead).
    $origC = @{}
    $origD = @{}
    for ($r = $start; $r -le $end; $r++) {
        $origC[$r] = $ws.Cells.Item($r, 3).Value2
        $origD[$r] = $ws.Cells.Item($r, 4).Value2
    }

    for ($r = $start; $r -lt $end; $r++) {
        # Shift: row r gets what row r+1 used to hold. The leading "'" keeps
        # Excel from re-parsing the date-shaped string back into a real date
        # (the source sheet stores Trimestre as plain text); re-applying the
        # "Normal" style afterwards drops the quote-prefix marker that the
        # apostrophe trick leaves behind, so the cell ends up styled exactly
        # like it started.
        $ws.Cells.Item($r, 3).Value = "'" + $origC[$r + 1]
        $ws.Cells.Item($r, 3).Style = "Normal"

        $nextD = $origD[$r + 1]
        if ($nextD -eq $null) {
            $ws.Cells.Item($r, 4).ClearContents()
        } else {
            $ws.Cells.Item($r, 4).Value = $nextD
        }
    }

    # Final row of the block: a new quarter appended to the series.
    $ws.Cells.Item($end, 3).Value = "'" + $block.newDate
    $ws.Cells.Item($end, 3).Style = "Normal"
    $ws.Cells.Item($end, 4).Value = $block.newValue
}
